$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New submission rows (8,9,10 in the "submission #" column -> rows 10,11,12) ---
# Row 10: submission #5, XGB, eta = 0.03
$ws.Range("B10").Value = 0.97869899999999999
$ws.Range("C6:C6").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = 42233
$ws.Range("D10").Value = "XGB"
$ws.Range("E10").Value = "eta = 0.03"

# Row 11: submission #6, XGB, eta = 0.015
$ws.Range("B11").Value = 0.97797199999999995
$ws.Range("C6:C6").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 42233
$ws.Range("D11").Value = "XGB"
$ws.Range("E11").Value = "eta = 0.015"

# Row 12: submission #7, XGB, eta=0.06
$ws.Range("B12").Value = 0.97828599999999999
$ws.Range("C6:C6").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = 42233
$ws.Range("D12").Value = "XGB"
$ws.Range("E12").Value = "eta=0.06"

# --- Best-fit column width refresh (approximation of Excel's autofit recalculation) ---
$ws.Columns.Item(1).ColumnWidth = 13.592447916666666
$ws.Columns.Item(3).ColumnWidth = 8.877604166666666
$ws.Columns.Item(4).ColumnWidth = 16.022135416666668
$ws.Columns.Item(5).ColumnWidth = 17.022135416666668
$ws.Columns.Item(6).ColumnWidth = 18.307291666666668
$ws.Columns.Item(7).ColumnWidth = 26.736979166666668

# --- Selection moved from H5 to G11 ---
[void]$ws.Range("G11").Select()
